$wb = $excel.ActiveWorkbook

# --- Rename the original "addVisitor" sheet to "addBrandNewVisitor" ---
$wsBrandNew = $wb.Worksheets.Item(1)
$wsBrandNew.Name = "addBrandNewVisitor"

# --- Insert the new "addVisitorForExistingIndividual" sheet right after it ---
# (Worksheets.Add(Before, After) places it between addBrandNewVisitor and sqlCount,
#  and the newly created sheet automatically becomes the active sheet.)
$wsExisting = $wb.Worksheets.Add($null, $wsBrandNew)
$wsExisting.Name = "addVisitorForExistingIndividual"

# --- Re-fetch the "sqlCount" sheet by name now that indices have shifted ---
$wsSqlCount = $wb.Worksheets.Item("sqlCount")

# --- Populate addBrandNewVisitor (sheet1) ---
$wsBrandNew.Range("A1").Value = "brandNewVisitorEmailId"
$wsBrandNew.Range("B1").Value = "lastName"
$wsBrandNew.Range("A2").Value = "SeleniumTest+v20191112113007@gmail.com"
$wsBrandNew.Range("B2").Value = "Test+v20191112113007"
$wsBrandNew.Range("A3").Value = "SeleniumTest+v20191112113325@gmail.com"
$wsBrandNew.Range("B3").Value = "Test+v20191112113325"

# Widen column A (source width 38.39 characters) and move the active
# selection to B1 on this sheet
$wsBrandNew.Columns.Item(1).ColumnWidth = 37.5
$wsBrandNew.Range("B1").Select()

# --- Populate addVisitorForExistingIndividual (new sheet) ---
$wsExisting.Range("A1").Value = "lastName"
$wsExisting.Range("B1").Value = "Test+v20191113110438"
$wsExisting.Range("A1").Select()

# --- Update sqlCount (now the 3rd sheet) ---
$wsSqlCount.Range("A1").Value = "sqlRecordCount"
$wsSqlCount.Range("B1").Value = "sqlColCount"

# "252" and "5" look numeric, so force them to stay text (matching the source
# workbook where they are stored as shared strings, not numbers): briefly mark
# the cell as Text, assign the value, then restore the General format so no
# extra style sticks to the cell.
$wsSqlCount.Range("A2").NumberFormat = "@"
$wsSqlCount.Range("A2").Value = "252"
$wsSqlCount.Range("A2").NumberFormat = "General"

$wsSqlCount.Range("B2").NumberFormat = "@"
$wsSqlCount.Range("B2").Value = "5"
$wsSqlCount.Range("B2").NumberFormat = "General"

# --- Make sure addVisitorForExistingIndividual ends up the active tab ---
$wsExisting.Activate()
